$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text content looks numeric (e.g. "7.1", "1.0") must be forced to
# text/string storage so Excel does not silently coerce them into numeric values.
$textCells = @("C16","L16","B21","K21","B22","K22","K23","B24","K24","B25","K25")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C16").Value = "7.1"
$ws.Range("L16").Value = "5.0"
$ws.Range("B21").Value = "1.0"
$ws.Range("K21").Value = "1.0"
$ws.Range("B22").Value = "1.0"
$ws.Range("K22").Value = "1.0"
$ws.Range("K23").Value = "1.0"
$ws.Range("B24").Value = "2.0"
$ws.Range("K24").Value = "1.0"
$ws.Range("B25").Value = "1.1"
$ws.Range("K25").Value = "1.0"

# Remaining cell updates (numeric values and plain text)
$ws.Range("B2").Value = 16
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = "Bowled"
$ws.Range("E2").Value = " Nuwan Pradeep"
$ws.Range("K2").Value = 71
$ws.Range("L2").Value = 20
$ws.Range("M2").Value = "NOT OUT"
$ws.Range("N2").Value = " "
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "Caught"
$ws.Range("E3").Value = " Chamika Karunarathne"
$ws.Range("K3").Value = 12
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = "LBW"
$ws.Range("N3").Value = " Imad Wasim"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "Bowled"
$ws.Range("E4").Value = " Maheesh Theekshana"
$ws.Range("M4").Value = "* NOT OUT"
$ws.Range("N4").Value = " "
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "Bowled"
$ws.Range("E5").Value = " Maheesh Theekshana"
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = " "
$ws.Range("N5").Value = " "
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "Bowled"
$ws.Range("E6").Value = " Maheesh Theekshana"
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = " "
$ws.Range("N6").Value = " "
$ws.Range("C7").Value = 15
$ws.Range("E7").Value = " Nuwan Pradeep"
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = " "
$ws.Range("N7").Value = " "
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "Bowled"
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = " "
$ws.Range("N8").Value = " "
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 3
$ws.Range("E9").Value = " Dushmantha Chameera"
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = " "
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "LBW"
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = " "
$ws.Range("N10").Value = " "
$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 3
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = " "
$ws.Range("N11").Value = " "
$ws.Range("B12").Value = 1
$ws.Range("D12").Value = "NOT OUT"
$ws.Range("E12").Value = " "
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = " "
$ws.Range("N12").Value = " "
$ws.Range("A16").Value = 91
$ws.Range("D16").Value = 43
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = 1
$ws.Range("M16").Value = 30
$ws.Range("A21").Value = "Chamika Karunarathne"
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 6
$ws.Range("J21").Value = "Haris Rauf"
$ws.Range("L21").Value = 22
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 22
$ws.Range("A22").Value = "Wanindu Hasaranga"
$ws.Range("C22").Value = 22
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 22
$ws.Range("J22").Value = "Shaheen Afridi"
$ws.Range("L22").Value = 24
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 24
$ws.Range("A23").Value = "Nuwan Pradeep"
$ws.Range("C23").Value = 29
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 14.5
$ws.Range("J23").Value = "Hasan Ali"
$ws.Range("L23").Value = 19
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 19
$ws.Range("A24").Value = "Maheesh Theekshana"
$ws.Range("C24").Value = 23
$ws.Range("E24").Value = 11.5
$ws.Range("J24").Value = "Imad Wasim"
$ws.Range("L24").Value = 12
$ws.Range("N24").Value = 12
$ws.Range("A25").Value = "Dushmantha Chameera"
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 10
$ws.Range("J25").Value = "Shadab Khan"
$ws.Range("L25").Value = 18
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 18

Write-Host "Applied match-report corrections."
